$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17
$ws.Cells.Item(17, 8).Value = 1303.2593
$ws.Cells.Item(17, 10).Value = 1475.619
$ws.Cells.Item(17, 12).Value = 4426.857
$ws.Cells.Item(17, 14).Value = -4762.857

# ALC row 40
$ws.Cells.Item(40, 8).Value = 10000
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()

# ALC row 43
$ws.Cells.Item(43, 8).Value = 1937.5
$ws.Cells.Item(43, 10).Value = 2000
$ws.Cells.Item(43, 12).Value = 2000
$ws.Cells.Item(43, 14).Value = -2138

# ALC row 70
$ws.Cells.Item(70, 8).Value = 1498.4445
$ws.Cells.Item(70, 10).Value = 1600
$ws.Cells.Item(70, 12).Value = 4800
$ws.Cells.Item(70, 14).Value = -5340

# ALC row 73
$ws.Cells.Item(73, 8).Value = 1498.4445
$ws.Cells.Item(73, 10).Value = 1600
$ws.Cells.Item(73, 12).Value = 4800
$ws.Cells.Item(73, 14).Value = -6672

# ALC row 141
$ws.Cells.Item(141, 8).Value = 3529.8
$ws.Cells.Item(141, 9).Value = 3529.8
$ws.Cells.Item(141, 11).Value = 10589.4
$ws.Cells.Item(141, 13).Value = -5409.400000000001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 17
$ws.Cells.Item(17, 8).Value = 8
$ws.Cells.Item(17, 9).Value = 8
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 8
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 14).ClearContents()
$ws.Cells.Item(17, 13).Value = 165

# ARM row 32
$ws.Cells.Item(32, 8).Value = 5207
$ws.Cells.Item(32, 9).Value = 5207
$ws.Cells.Item(32, 11).Value = 5207
$ws.Cells.Item(32, 13).Value = -4920

# ARM row 63
$ws.Cells.Item(63, 8).Value = 6163.9
$ws.Cells.Item(63, 10).Value = 10872.5
$ws.Cells.Item(63, 12).Value = 10872.5
$ws.Cells.Item(63, 14).Value = -12244.5

# ARM row 66
$ws.Cells.Item(66, 8).Value = 6163.9
$ws.Cells.Item(66, 10).Value = 10872.5
$ws.Cells.Item(66, 12).Value = 54362.5
$ws.Cells.Item(66, 14).Value = -61226.5

# ARM row 102
$ws.Cells.Item(102, 8).Value = 2237.5
$ws.Cells.Item(102, 9).Value = 2237.5
$ws.Cells.Item(102, 11).Value = 2237.5
$ws.Cells.Item(102, 13).Value = -615.5

# ARM row 132
$ws.Cells.Item(132, 8).Value = 3967.25
$ws.Cells.Item(132, 9).Value = 3770.7273
$ws.Cells.Item(132, 11).Value = 11312.1819
$ws.Cells.Item(132, 13).Value = -8782.1819

$ws = $wb.Worksheets.Item("BSM")
# BSM row 131
$ws.Cells.Item(131, 8).Value = 88900
$ws.Cells.Item(131, 10).Value = 88900
$ws.Cells.Item(131, 12).Value = 88900
$ws.Cells.Item(131, 14).Value = -98980

# BSM row 134
$ws.Cells.Item(134, 8).Value = 2269.5715
$ws.Cells.Item(134, 9).Value = 2577.4
$ws.Cells.Item(134, 11).Value = 7732.200000000001
$ws.Cells.Item(134, 13).Value = -5197.200000000001

# BSM row 141
$ws.Cells.Item(141, 8).Value = 159979
$ws.Cells.Item(141, 9).Value = 159970
$ws.Cells.Item(141, 11).Value = 159970
$ws.Cells.Item(141, 13).Value = -154790

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 2272.3076
$ws.Cells.Item(31, 9).Value = 2101.5715
$ws.Cells.Item(31, 11).Value = 2101.5715
$ws.Cells.Item(31, 13).Value = -1806.5715

# CRP row 34
$ws.Cells.Item(34, 8).Value = 2272.3076
$ws.Cells.Item(34, 9).Value = 2101.5715
$ws.Cells.Item(34, 11).Value = 2101.5715
$ws.Cells.Item(34, 13).Value = -1899.5715

# CRP row 58
$ws.Cells.Item(58, 8).Value = 1729.5161
$ws.Cells.Item(58, 9).Value = 1809.8334
$ws.Cells.Item(58, 10).Value = 1454.1428
$ws.Cells.Item(58, 11).Value = 1809.8334
$ws.Cells.Item(58, 12).Value = 1454.1428
$ws.Cells.Item(58, 13).Value = -1606.8334
$ws.Cells.Item(58, 14).Value = -1860.1428

# CRP row 107
$ws.Cells.Item(107, 8).Value = 101464.6
$ws.Cells.Item(107, 9).Value = 144378
$ws.Cells.Item(107, 10).Value = 1333.3334
$ws.Cells.Item(107, 11).Value = 144378
$ws.Cells.Item(107, 12).Value = 1333.3334
$ws.Cells.Item(107, 13).Value = -142458
$ws.Cells.Item(107, 14).Value = -5173.3334

# CRP row 122
$ws.Cells.Item(122, 8).Value = 1406.7778
$ws.Cells.Item(122, 9).Value = 1457.625
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 4372.875
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -1922.875
$ws.Cells.Item(122, 14).Value = -7900

# CRP row 132
$ws.Cells.Item(132, 8).Value = 127924.75
$ws.Cells.Item(132, 9).Value = 168566.67
$ws.Cells.Item(132, 11).Value = 505700.01
$ws.Cells.Item(132, 13).Value = -503170.01

# CRP row 134
$ws.Cells.Item(134, 8).Value = 4944
$ws.Cells.Item(134, 9).Value = 5062.25
$ws.Cells.Item(134, 11).Value = 15186.75
$ws.Cells.Item(134, 13).Value = -12651.75

# CRP row 136
$ws.Cells.Item(136, 8).Value = 1729.5161
$ws.Cells.Item(136, 9).Value = 1809.8334
$ws.Cells.Item(136, 10).Value = 1454.1428
$ws.Cells.Item(136, 11).Value = 5429.5002
$ws.Cells.Item(136, 12).Value = 4362.428400000001
$ws.Cells.Item(136, 13).Value = -2879.5002
$ws.Cells.Item(136, 14).Value = -9462.4284

$ws = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws.Cells.Item(11, 8).Value = 825000
$ws.Cells.Item(11, 9).Value = 900000
$ws.Cells.Item(11, 10).Value = 600000
$ws.Cells.Item(11, 11).Value = 2700000
$ws.Cells.Item(11, 12).Value = 1800000
$ws.Cells.Item(11, 13).Value = -2699860
$ws.Cells.Item(11, 14).Value = -1800280

# CUL row 14
$ws.Cells.Item(14, 8).Value = 345
$ws.Cells.Item(14, 9).Value = 345
$ws.Cells.Item(14, 11).Value = 1035
$ws.Cells.Item(14, 13).Value = -862

# CUL row 131
$ws.Cells.Item(131, 8).Value = 1923.9
$ws.Cells.Item(131, 9).Value = 2413.3333
$ws.Cells.Item(131, 11).Value = 7239.999899999999
$ws.Cells.Item(131, 13).Value = -2199.999899999999

# CUL row 134
$ws.Cells.Item(134, 8).Value = 1899.6666
$ws.Cells.Item(134, 9).Value = 1899.6666
$ws.Cells.Item(134, 11).Value = 5698.9998
$ws.Cells.Item(134, 13).Value = -628.9997999999996

# CUL row 138
$ws.Cells.Item(138, 8).Value = 2526.2
$ws.Cells.Item(138, 9).Value = 500
$ws.Cells.Item(138, 10).Value = 3032.75
$ws.Cells.Item(138, 11).Value = 1500
$ws.Cells.Item(138, 12).Value = 9098.25
$ws.Cells.Item(138, 14).Value = -19378.25
$ws.Cells.Item(138, 13).Value = 3640

# CUL row 139
$ws.Cells.Item(139, 8).Value = 202415.4
$ws.Cells.Item(139, 9).Value = 334670.34
$ws.Cells.Item(139, 10).Value = 4033
$ws.Cells.Item(139, 11).Value = 1004011.02
$ws.Cells.Item(139, 12).Value = 12099
$ws.Cells.Item(139, 13).Value = -998871.02
$ws.Cells.Item(139, 14).Value = -22379

$ws = $wb.Worksheets.Item("GSM")
# GSM row 113
$ws.Cells.Item(113, 8).Value = 2287.75
$ws.Cells.Item(113, 9).Value = 2541.4
$ws.Cells.Item(113, 10).Value = 1865
$ws.Cells.Item(113, 11).Value = 2541.4
$ws.Cells.Item(113, 12).Value = 1865
$ws.Cells.Item(113, 13).Value = -371.4000000000001
$ws.Cells.Item(113, 14).Value = -6205

# GSM row 126
$ws.Cells.Item(126, 8).Value = 3999.5
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 13).Value = -9530

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Cells.Item(40, 8).Value = 4400
$ws.Cells.Item(40, 9).Value = 4400
$ws.Cells.Item(40, 11).Value = 4400
$ws.Cells.Item(40, 13).Value = -4264

# LTW row 68
$ws.Cells.Item(68, 8).Value = 2514.6667
$ws.Cells.Item(68, 9).Value = 2514.6667
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 2514.6667
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -1765.6667
$ws.Cells.Item(68, 14).ClearContents()

# LTW row 71
$ws.Cells.Item(71, 8).Value = 2514.6667
$ws.Cells.Item(71, 9).Value = 2514.6667
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 12573.3335
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = -8829.3335
$ws.Cells.Item(71, 14).ClearContents()
